# This script applies the "Added handling of common packages" commit to the
# mall-gateway_structure workbook.
#
# Summary of the semantic change (derived from the OOXML diff):
#
# 1) Sheet "methodNumberOfLines": four no-arg constructor rows are removed
#    (the structure-extraction tool now skips/merges these, e.g. because of
#    the new "common packages" handling):
#       com.macro.mall.config.IgnoreUrlsConfig        -> IgnoreUrlsConfig()
#       com.macro.mall.config.GlobalCorsConfig        -> GlobalCorsConfig()
#       com.macro.mall.MallGatewayApplication         -> MallGatewayApplication()
#       com.macro.mall.filter.IgnoreUrlsRemoveJwtFilter -> IgnoreUrlsRemoveJwtFilter()
#    All other rows keep their relative order (dimension shrinks from
#    A1:C41 to A1:C37).
#
# 2) Sheet "classFields": the fields that belong to
#    com.macro.mall.config.ResourceServerConfig (rows 10-14) are reordered to
#       authorizationManager, ignoreUrlsConfig, ignoreUrlsRemoveJwtFilter,
#       restfulAccessDeniedHandler, restAuthenticationEntryPoint
#    and the fields of com.macro.mall.authorization.AuthorizationManager
#    (rows 15-16) are reordered to
#       ignoreUrlsConfig, redisTemplate

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) methodNumberOfLines: delete the four no-arg constructor rows.
#    Delete from the bottom up so earlier row numbers stay valid.
# ---------------------------------------------------------------------------
$wsMethods = $wb.Worksheets.Item("methodNumberOfLines")

$wsMethods.Rows.Item(38).Delete()  # com.macro.mall.filter.IgnoreUrlsRemoveJwtFilter | IgnoreUrlsRemoveJwtFilter()
$wsMethods.Rows.Item(30).Delete()  # com.macro.mall.MallGatewayApplication | MallGatewayApplication()
$wsMethods.Rows.Item(28).Delete()  # com.macro.mall.config.GlobalCorsConfig | GlobalCorsConfig()
$wsMethods.Rows.Item(11).Delete()  # com.macro.mall.config.IgnoreUrlsConfig | IgnoreUrlsConfig()

# ---------------------------------------------------------------------------
# 2) classFields: reorder the field rows for ResourceServerConfig and
#    AuthorizationManager. Column A (Class Name) and C (Field Modifier) are
#    unchanged for these rows, only B (Field Name) / D (Field Type) move.
# ---------------------------------------------------------------------------
$wsFields = $wb.Worksheets.Item("classFields")

# com.macro.mall.config.ResourceServerConfig fields (rows 10-14)
$wsFields.Range("B10").Value = "authorizationManager"
$wsFields.Range("D10").Value = "com.macro.mall.authorization.AuthorizationManager"

$wsFields.Range("B11").Value = "ignoreUrlsConfig"
$wsFields.Range("D11").Value = "com.macro.mall.config.IgnoreUrlsConfig"

$wsFields.Range("B12").Value = "ignoreUrlsRemoveJwtFilter"
$wsFields.Range("D12").Value = "com.macro.mall.filter.IgnoreUrlsRemoveJwtFilter"

$wsFields.Range("B13").Value = "restfulAccessDeniedHandler"
$wsFields.Range("D13").Value = "com.macro.mall.component.RestfulAccessDeniedHandler"

$wsFields.Range("B14").Value = "restAuthenticationEntryPoint"
$wsFields.Range("D14").Value = "com.macro.mall.component.RestAuthenticationEntryPoint"

# com.macro.mall.authorization.AuthorizationManager fields (rows 15-16)
$wsFields.Range("B15").Value = "ignoreUrlsConfig"
$wsFields.Range("D15").Value = "com.macro.mall.config.IgnoreUrlsConfig"

$wsFields.Range("B16").Value = "redisTemplate"
$wsFields.Range("D16").Value = "org.springframework.data.redis.core.RedisTemplate"

# row 17 (com.macro.mall.filter.IgnoreUrlsRemoveJwtFilter -> ignoreUrlsConfig)
# is unchanged in content, left as-is.

Write-Host "Done. methodNumberOfLines rows: $($wsMethods.UsedRange.Rows.Count); classFields rows: $($wsFields.UsedRange.Rows.Count)"
